$d = $word.ActiveDocument

$pairs = @(
    @("675÷9=", "497÷4="),
    @("761÷9=", "839÷8="),
    @("300÷6=", "246÷3="),
    @("966÷6=", "663÷5="),
    @("244÷2=", "332÷4="),
    @("776÷6=", "820÷2="),
    @("504÷4=", "597÷8="),
    @("337÷7=", "983÷3="),
    @("345÷6=", "443÷5="),
    @("553÷9=", "409÷8="),
    @("707÷6=", "255÷2="),
    @("356÷9=", "851÷3="),
    @("736÷6=", "223÷5="),
    @("889÷2=", "209÷2="),
    @("597÷6=", "455÷3="),
    @("786÷8=", "301÷2="),
    @("701÷2=", "517÷4="),
    @("637÷6=", "741÷7="),
    @("650÷4=", "387÷6="),
    @("103÷8=", "319÷7="),
    @("623÷4=", "326÷5="),
    @("698÷4=", "134÷8="),
    @("525÷3=", "305÷8="),
    @("116÷5=", "193÷6="),
    @("327÷2=", "481÷4=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
